$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column (C) for rows 2 through 28
# from Excel serial date 45516 (2024-08-12) to 45517 (2024-08-13).
for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45516) {
        $cell.Value2 = 45517
    }
}
